# Applies the changes described by the commit diff:
# 1. Rename worksheet "Sheet1" -> "Apple Cinema"
# 2. Add a new row of data to the "Contact" sheet (row 5), make it the active sheet/cell
# 3. Make "Contact" the active sheet (activeTab changes from 5 to 2), so it is no longer
#    on "Sheet1"/"Apple Cinema"

$wb = $excel.ActiveWorkbook

# 1. Rename "Sheet1" to "Apple Cinema"
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "Apple Cinema"

# 2. Append new row of contact-form data onto the "Contact" sheet
$contact = $wb.Worksheets.Item("Contact")
$contact.Range("A5").Value = "Tester"
$contact.Range("B5").Value = "deekshavishwakarma@yahoo.com"
$contact.Range("C5").Value = "deeksha"

# Select B5 on the Contact sheet and make it the active sheet
$contact.Range("B5").Select()
$contact.Activate()
